$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update password column (B2:B11) from 12345 to 12345678
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Value = 12345678
}

# Update the selected cell/range
$ws.Range("B14").Select()
